$d = $word.ActiveDocument
$d.Content.Find.Execute("LICENCE RENEWAL 2021 – 2022", $false, $false, $false, $false, $false, $true, 1, $false, "LICENCE RENEWAL 2022 – 2023", 2)
